$d = $word.ActiveDocument

# Locate the existing backlog item "Manage flowchart execution speed" -
# the new items are added right after it, before the trailing blank
# "List Paragraph" entries at the end of the backlog list.
$anchor = $d.Content
$found = $anchor.Find.Execute("Manage flowchart execution speed", $false, $false, `
    $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not locate anchor paragraph 'Manage flowchart execution speed'"
}

# Move past the end of the found text's paragraph mark so the insertion
# point sits right between that paragraph and the next one.
$anchor.Collapse(0)
[void]$anchor.MoveEnd(1, 1)
$anchor.Collapse(0)

# Build a brand new Range object at that position (re-using the Find
# result range directly can confuse InsertXML into overwriting the
# paragraph that was just matched).
$insertionPoint = $d.Range($anchor.Start, $anchor.End)

$w = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

# New bulleted backlog item, using the same numbering/list style as the
# surrounding items ("List Paragraph" / numId 2), plus a trailing blank
# "List Paragraph" line (no numbering) like the ones further down.
$newItemPPr = "<w:pPr><w:pStyle w:val='ListParagraph'/><w:numPr><w:ilvl w:val='0'/><w:numId w:val='2'/></w:numPr><w:jc w:val='both'/><w:rPr><w:rFonts w:ascii='Verdana' w:hAnsi='Verdana'/><w:sz w:val='24'/></w:rPr></w:pPr>"
$newItemRPr = "<w:rPr><w:rFonts w:ascii='Verdana' w:hAnsi='Verdana'/><w:sz w:val='24'/></w:rPr>"
$newItemRun = "<w:r>$newItemRPr<w:t>Create different user themes</w:t></w:r>"
$newItemPara = "<w:p $w>$newItemPPr$newItemRun</w:p>"

$blankPPr = "<w:pPr><w:pStyle w:val='ListParagraph'/><w:jc w:val='both'/><w:rPr><w:rFonts w:ascii='Verdana' w:hAnsi='Verdana'/><w:sz w:val='24'/></w:rPr></w:pPr>"
$blankPara = "<w:p $w>$blankPPr</w:p>"

[void]$insertionPoint.InsertXML($newItemPara + $blankPara)
